$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.654535055160522
$ws.Range("B1").Value = 2.201472282409668
$ws.Range("C1").Value = 4.29674243927002
$ws.Range("D1").Value = 4.619184494018555
$ws.Range("E1").Value = 1.815486431121826
